$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Sexo" labels used throughout the data: Masculino -> Hombre, Femenino -> Mujer
$ws.Cells.Replace("Masculino", "Hombre")
$ws.Cells.Replace("Femenino", "Mujer")

# Move the active selection to T19 (matches the author's last cursor position)
$ws.Activate()
$ws.Range("T19").Select()
